$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to be stored as Text (matches
# how the source data was authored as inline strings), then restore the
# cells NumberFormat back to General so no stray formatting diff remains.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '29.123.94'
$ws.Cells.Item(2, 5).Value = '  -1.77%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '1.836.12'
$ws.Cells.Item(3, 5).Value = '  -1.41%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '0.9993'
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5, 5).Value = '  -2.41%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '0.6806'
$ws.Cells.Item(6, 5).Value = '  -2.72%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '1.0000'
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '0.2988'
$ws.Cells.Item(8, 5).Value = '  -2.72%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.07439'
$ws.Cells.Item(9, 5).Value = '  -3.89%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -2.26%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.07648'
$ws.Cells.Item(11, 5).Value = '  -1.62%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '1.838.51'
$ws.Cells.Item(12, 5).Value = '  -1.30%  '

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) '5.025'

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '0.6785'
$ws.Cells.Item(14, 5).Value = '  -2.10%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '86.83'
$ws.Cells.Item(15, 5).Value = '  -6.01%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '6.154'
$ws.Cells.Item(16, 5).Value = '  -6.33%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '29.107.37'

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '0.000008219'
$ws.Cells.Item(18, 5).Value = '  -1.84%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '229.24'
$ws.Cells.Item(19, 5).Value = '  -5.29%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -2.26%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.04%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '7.336'
$ws.Cells.Item(22, 5).Value = '  -3.75%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '0.9998'
$ws.Cells.Item(23, 5).Value = '  -0.03%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '161.24'
$ws.Cells.Item(24, 5).Value = '  +1.03%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '0.1430'
$ws.Cells.Item(25, 5).Value = '  -5.32%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '8.709'
$ws.Cells.Item(26, 5).Value = '  -2.41%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -1.62%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -2.53%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '4.246'
$ws.Cells.Item(29, 5).Value = '  -0.39%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '4.139'
$ws.Cells.Item(30, 5).Value = '  -1.41%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '1.189'
$ws.Cells.Item(31, 5).Value = '  -0.42%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '0.05339'
$ws.Cells.Item(32, 5).Value = '  +4.44%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '0.7540'
$ws.Cells.Item(33, 5).Value = '  -3.95%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '1.846'
$ws.Cells.Item(34, 5).Value = '  -3.01%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -2.42%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '2.682'
$ws.Cells.Item(36, 5).Value = '  -0.16%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '1.314.08'
$ws.Cells.Item(37, 5).Value = '  -1.46%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '0.01821'
$ws.Cells.Item(38, 5).Value = '  -3.24%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '2.713'
$ws.Cells.Item(39, 5).Value = '  -0.98%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '0.9363'
$ws.Cells.Item(40, 5).Value = '  -3.26%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '6.064'
$ws.Cells.Item(41, 5).Value = '  +1.62%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '104.92'
$ws.Cells.Item(42, 5).Value = '  -1.52%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '0.08318'
$ws.Cells.Item(43, 5).Value = '  +30.87%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '0.9989'
$ws.Cells.Item(44, 5).Value = '  -0.07%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '1.981.54'
$ws.Cells.Item(45, 5).Value = '  -1.38%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '0.5181'
$ws.Cells.Item(46, 5).Value = '  -0.66%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Cells.Item(47, 4) '0.00000000121'
$ws.Cells.Item(47, 5).Value = '  -3.92%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Aave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Cells.Item(48, 4) '64.06'
$ws.Cells.Item(48, 5).Value = '  -1.29%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '1.768'
$ws.Cells.Item(49, 5).Value = '  -1.24%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '9.376'
$ws.Cells.Item(50, 5).Value = '  -4.07%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '0.05931'
$ws.Cells.Item(51, 5).Value = '  +0.33%  '
